$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.276191830635071
$ws.Range("B1").Value = 2.247664928436279
$ws.Range("C1").Value = 4.309806823730469
$ws.Range("D1").Value = 3.008450269699097
$ws.Range("E1").Value = 1.352496385574341
